# Add data for 2021-12-17: update sheet name, row 13 (November) 2021 totals,
# row 14 (December through 12-xx) values, and row 15 (Total) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet / sheet tab title to reflect the new "through" date.
$ws.Name = "Through 2021-12-09"

# Row 13 (November): update 2021 totals (no_arrest_made / arrest_rate)
$ws.Range("U13").Value = 196
$ws.Range("V13").Value = 0.0249

# Row 14 (December through 12-09): label + per-year updates
$ws.Range("A14").Value = "December (through 12-09)"

$ws.Range("F14").Value = 24
$ws.Range("G14").Value = 0.0769

$ws.Range("I14").Value = 33
$ws.Range("J14").Value = 0.0833

$ws.Range("L14").Value = 19
$ws.Range("M14").Value = 0.0952

$ws.Range("O14").Value = 10
$ws.Range("P14").Value = 0.2308

$ws.Range("R14").Value = 40
$ws.Range("S14").Value = 0.0476

$ws.Range("U14").Value = 71

# Row 15 (Total): update per-year totals
$ws.Range("F15").Value = 527
$ws.Range("G15").Value = 0.1053

$ws.Range("I15").Value = 791
$ws.Range("J15").Value = 0.077

$ws.Range("L15").Value = 627
$ws.Range("M15").Value = 0.1081

$ws.Range("O15").Value = 490
$ws.Range("P15").Value = 0.1042

$ws.Range("R15").Value = 1240
$ws.Range("S15").Value = 0.0505

$ws.Range("U15").Value = 1615
$ws.Range("V15").Value = 0.0578
